# Saldo.xlsx update: refresh balances for three existing accounts and
# add/move the BLUEMETRIX (001761119) row to its new sorted position
# (the sheet is kept sorted descending by the "Saldo" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update balances that changed in place (rows 6, 7, 8) -----------------
$ws.Cells.Item(6, 3).Value = 102114.63
$ws.Cells.Item(7, 3).Value = 100000
$ws.Cells.Item(8, 3).Value = 46399.01

# --- Move BLUEMETRIX (001761119) from row 232 to its new sorted slot ------
# Its balance grew from 41.64 to 204.87, so it now belongs right before
# NASSIM / 004448303 (row 133), between CLAUDIA / 004693697 (213.07) and
# NASSIM (204.44).
$ws.Rows(133).Insert()
$ws.Range("A233:C233").Cut($ws.Range("A133:C133"))
$ws.Rows(233).Delete()
$ws.Cells.Item(133, 3).Value = 204.87
